# Make plot and visit forms more consistent:
# rename the "name" field of the plot survey form to "plot_name",
# and update the settings sheet's instance_name reference to match.

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$choices  = $wb.Worksheets.Item("choices")
$settings = $wb.Worksheets.Item("settings")

# survey!D3 held the literal "name" (the survey row whose `name` column
# says the field is called "name") -> rename it to "plot_name"
$survey.Range("D3").Value = "plot_name"

# settings!B5 is the instance_name setting, which referenced the field
# called "name" -> it now needs to reference "plot_name" too
$settings.Range("B5").Value = "plot_name"

# Update the remembered selections on each sheet
$survey.Range("D24").Select()
$choices.Range("C6").Select()
$settings.Range("B5").Select()

# settings is now the active/selected tab
$settings.Activate()
